$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 17 - separator row: only A17 has text, whole row gets a thin
# bottom border (A/B use the "plain" style, C/D/E use the "font1" style)
# -----------------------------------------------------------------
$ws.Range("A17").Value = "SCRIPT/P02P01A/um2102.ssb"

$row17 = $ws.Range("A17:E17")
$row17.RowHeight = 43.2
$bottom17 = $row17.Borders.Item(9)
$bottom17.LineStyle = 1
$bottom17.Weight = 2

$ws.Range("C17:E17").Font.Size = 8

# -----------------------------------------------------------------
# Row 18 - header-like row with thin borders on top AND bottom
# -----------------------------------------------------------------
$ws.Range("A18").Value = "SCRIPT/T01P01A/um2201.ssb"
$ws.Range("B18").Value = 63
$ws.Range("C18").Value = " The whole world\'s time is about\nto stop?!"
$ws.Range("D18").Value = " Скоро во всём мире остановится\nвремя?!"
$ws.Range("E18").Value = " Òëïñï âï âòæí íéñå ïòóàîïâéóòÿ\nâñåíÿ?!"

$row18 = $ws.Range("A18:E18")
$row18.RowHeight = 43.2
$top18 = $row18.Borders.Item(8)
$top18.LineStyle = 1
$top18.Weight = 2
$bottom18 = $row18.Borders.Item(9)
$bottom18.LineStyle = 1
$bottom18.Weight = 2

$ws.Range("C18:E18").Font.Size = 8

# -----------------------------------------------------------------
# Row 19 - regular data row (same style family as rows 2-16)
# -----------------------------------------------------------------
$ws.Range("A19").Value = "SCRIPT/P02P01A/um2402.ssb"
$ws.Range("B19").Value = 38
$ws.Range("C19").Value = " The world\'s going to be ruined if\nsomething\'s not done, is that right?"
$ws.Range("D19").Value = " Миру придёт конец если ничего\nне сделать, так?"
$ws.Range("E19").Value = " Íéñô ðñéäæó ëïîåø åòìé îéœåãï\nîå òäåìàóû, óàë?"
$ws.Range("A19:E19").RowHeight = 43.2
$ws.Range("C19:E19").Font.Size = 8

# -----------------------------------------------------------------
# Row 20 - regular data row
# -----------------------------------------------------------------
$ws.Range("A20").Value = "SCRIPT/P02P01A/um2502.ssb"
$ws.Range("B20").Value = 41
$ws.Range("C20").Value = " But what can I do to stop it?\nI can only cheer you on while lolling in the\nHot Spring…"
$ws.Range("D20").Value = " Но что я могу сделать, чтобы\nэтого не случилось? Только болеть за вас,\nотмокая в Горячих Источниках..."
$ws.Range("E20").Value = " Îï œóï ÿ íïãô òäåìàóû, œóïáú\nüóïãï îå òìôœéìïòû? Óïìûëï áïìåóû èà âàò,\nïóíïëàÿ â Ãïñÿœéö Éòóïœîéëàö..."
$ws.Range("A20:E20").RowHeight = 43.2
$ws.Range("C20:E20").Font.Size = 8

# -----------------------------------------------------------------
# Row 21 - regular row but with no value in column A and default row height
# -----------------------------------------------------------------
$ws.Range("B21").Value = 44
$ws.Range("C21").Value = " But go do it![K] Save the world!"
$ws.Range("D21").Value = " Но, вперёд![K] Спасите мир!"
$ws.Range("E21").Value = " Îï, âðåñæä![K] Òðàòéóå íéñ!"
$ws.Range("C21:E21").Font.Size = 8

# -----------------------------------------------------------------
# View state: scrolled so row 16 is at top, D19 is the active cell
# -----------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("D19").Select()
